# This script applies the commit "[ADDITIONAL SCRAPING] added code to scrape
# more data about a player's batting performance in a match, also updated
# the excel sheets" to the workbook.
#
# Summary of changes:
#  1. A new worksheet "Player Info" is inserted as the first sheet, holding
#     basic player bio data (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  2. On the existing "ODI Batting" sheet, column D (MATCH_CARD_LINK) is
#     renamed to MATCH_CODE and its values are trimmed from full scorecard
#     URLs down to just the numeric match code.
#  3. On the existing "ODI Bowling" sheet, column B (MATCH_CARD_LINK) is
#     renamed to MATCH_CODE and its values are likewise trimmed down to the
#     numeric match code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet in front of "ODI Batting".
#    NOTE: worksheet references are positional, so grab the other
#    sheets again (by name) *after* the insert has shifted indices.
# ---------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$infoSheet.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------
# Fill in the new "Player Info" sheet
# ---------------------------------------------------------------------
$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header formatting used by the
# other sheets' header rows.
$infoSheet.Range("A1:D1").Font.Bold = $true
$infoSheet.Range("A1:D1").Borders.LineStyle = 1
$infoSheet.Range("A1:D1").HorizontalAlignment = -4108
$infoSheet.Range("A1:D1").VerticalAlignment = -4160

$infoSheet.Range("A2:D2").NumberFormat = "@"
$infoSheet.Range("A2").Value = "6655"
$infoSheet.Range("B2").Value = "Dunith Nethmika Wellalage"
$infoSheet.Range("C2").Value = "Left Handed"
$infoSheet.Range("D2").Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------
# 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
# ---------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2 = "4594"
    3 = "4597"
    4 = "4600"
    5 = "4601"
    6 = "4603"
    7 = "4675"
    8 = "4687"
    9 = "4689"
    10 = "4691"
}

$battingSheet.Range("D2:D10").NumberFormat = "@"
foreach ($row in $battingCodes.Keys) {
    $battingSheet.Cells.Item($row, 4).Value = $battingCodes[$row]
}

# ---------------------------------------------------------------------
# 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
# ---------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4594"
    3 = "4597"
    4 = "4600"
    5 = "4601"
    6 = "4603"
    7 = "4675"
    8 = "4687"
    9 = "4689"
}

$bowlingSheet.Range("B2:B9").NumberFormat = "@"
foreach ($row in $bowlingCodes.Keys) {
    $bowlingSheet.Cells.Item($row, 2).Value = $bowlingCodes[$row]
}

Write-Host "Sheets after edit:"
foreach ($ws in $wb.Worksheets) {
    Write-Host " -" $ws.Name
}
